$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to Text format so that numeric-looking
# values (e.g. "1.006") are stored as text strings, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Cell updates per the refreshed crypto price/volume snapshot ---
$ws.Range("D2").Value = "26.152.78"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.654.98"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "218.54"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "0.5190"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "0.2669"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "0.06311"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "21.05"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.654.09"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "4.427"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "1.883.69"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "0.5460"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "0.0₅8223"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "64.84"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").Value = "26.215.44"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("D20").Value = "4.668"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "192.06"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "6.106"
$ws.Range("E23").Value = "  -4.79%  "
$ws.Range("D24").Value = "1.009"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "137.30"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("D26").Value = "0.1234"
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("D27").Value = "7.243"
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").Value = "16.10"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").Value = "1.413"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "0.06013"
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("D31").Value = "1.286"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "3.548"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "3.335"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").Value = "1.645"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").Value = "0.9796"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "2.779"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "0.5922"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").Value = "5.958"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "0.8669"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "1.038.24"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("D44").Value = "99.74"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "1.798.25"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "57.12"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈108"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "8.111"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "0.05182"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "1.474"

# Remove the temporary text-number-format style so the cells end up with
# no explicit style index (matching the original workbook formatting).
$ws.Range("D2:D51").Style = "Normal"

